$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 80 — this shifts the existing rows 80-183
# (and their formatting) down to 81-184, matching the diff's "everything
# slides down by one, new dimension A1:R184" shape.
$ws.Rows.Item(80).Insert()

# Populate the newly-inserted row 80 with the new data record (same
# constant fields as the rest of the Cilantro/Coquimbo block, new date
# and new Volumen figure; Precio minimo/maximo/promedio/$/Kg match what
# had previously been in this slot).
$ws.Cells.Item(80, 1).Value = 8
$ws.Cells.Item(80, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(80, 3).Value = "Coquimbo"
$ws.Cells.Item(80, 4).Value = 44894
$ws.Cells.Item(80, 5).Value = 4
$ws.Cells.Item(80, 6).Value = 100112040
$ws.Cells.Item(80, 7).Value = "Cilantro"
$ws.Cells.Item(80, 8).Value = "Sin especificar"
$ws.Cells.Item(80, 9).Value = "Primera"
$ws.Cells.Item(80, 10).Value = 2400
$ws.Cells.Item(80, 11).Value = 1500
$ws.Cells.Item(80, 12).Value = 2000
$ws.Cells.Item(80, 13).Value = 1750
$ws.Cells.Item(80, 14).Value = "`$/atado 1 a 1,5 kilos"
$ws.Cells.Item(80, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(80, 16).Value = 1167
$ws.Cells.Item(80, 17).Value = 1.5
$ws.Cells.Item(80, 18).Value = "Hortaliza"
